$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.988.69"
$ws.Range("E2").Value = "  +0.23%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.352.53"
$ws.Range("E3").Value = "  +0.04%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("E5").Value = "  +0.63%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "238.66"
$ws.Range("E6").Value = "  +1.18%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "73.55"
$ws.Range("E7").Value = "  +0.88%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.581"
$ws.Range("E9").Value = "  +6.63%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0998"
$ws.Range("E10").Value = "  +0.82%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "57.18"
$ws.Range("E11").Value = "  +0.18%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "31.85"
$ws.Range("E12").Value = "  +13.25%  "
$ws.Range("E13").Value = "  +0.81%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.20"
$ws.Range("E14").Value = "  +8.00%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.703.33"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "16.56"
$ws.Range("E16").Value = "  -0.80%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.898"
$ws.Range("E17").Value = "  +0.86%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.351.52"
$ws.Range("E18").Value = "  +0.54%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "43.860.05"
$ws.Range("E19").Value = "  -0.04%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0000101"
$ws.Range("E20").Value = "  -0.26%  "
$ws.Range("E21").Value = "  +4.25%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "76.67"
$ws.Range("E22").Value = "  -0.74%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "255.34"
$ws.Range("E23").Value = "  +0.53%  "
$ws.Range("E24").Value = "  +21.24%  "
$ws.Range("E25").Value = "  -0.06%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.68"
$ws.Range("E26").Value = "  -1.81%  "
$ws.Range("E27").Value = "  -0.50%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.65"
$ws.Range("E28").Value = "  +0.75%  "
$ws.Range("E29").Value = "  -0.61%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "22.62"
$ws.Range("E30").Value = "  +0.94%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "175.42"
$ws.Range("E31").Value = "  +1.58%  "
$ws.Range("E32").Value = "  -3.11%  "
$ws.Range("E33").Value = "  +3.10%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0754"
$ws.Range("E34").Value = "  +5.91%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.25"
$ws.Range("E35").Value = "  +1.73%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.35"
$ws.Range("E36").Value = "  +3.83%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.71"
$ws.Range("E37").Value = "  -5.79%  "
$ws.Range("B38").Value = "THORChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.32"
$ws.Range("E38").Value = "  -1.66%  "
$ws.Range("B39").Value = "LidoDAOToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.35"
$ws.Range("E39").Value = "  -2.93%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0279"
$ws.Range("E40").Value = "  +4.48%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.111"
$ws.Range("E41").Value = "  +13.76%  "
$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "9.14"
$ws.Range("E42").Value = "  +3.52%  "
$ws.Range("B43").Value = "InjectiveProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "19.13"
$ws.Range("E43").Value = "  -2.36%  "
$ws.Range("B44").Value = "Algorand"
$ws.Range("C44").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.201"
$ws.Range("E44").Value = "  +10.44%  "
$ws.Range("E45").Value = "  +0.00%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.70"
$ws.Range("E46").Value = "  +5.77%  "
$ws.Range("B47").Value = "MultiversX"
$ws.Range("C47").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "57.36"
$ws.Range("E47").Value = "  +9.38%  "
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.47"
$ws.Range("E48").Value = "  +7.82%  "
$ws.Range("B49").Value = "TrustWalletToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.24"
$ws.Range("E49").Value = "  +0.30%  "
$ws.Range("E50").Value = "  +0.46%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "99.52"
$ws.Range("E51").Value = "  +2.38%  "
